# Edit VIRGINIA_2019.xlsx per commit "Fixing network data cleanining scripts"
#  1. Drop the trailing metadata/footer rows (sample size, source, etc.)
#  2. Rename the header row to short machine-friendly column names
#  3. Title-case the Spanish connector words ("de", "del", "el", "la",
#     "los", "las", "y") inside state/municipality names
#  4. A handful of D-column percentages shift by a single ULP because the
#     upstream pipeline recomputed them - pin those exact cells too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the trailing metadata rows (now rows 1088-1092) -----------
$ws.Rows("1088:1092").Delete()

# --- 2. Rename header row -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 3. Title-case Spanish connector words in state/municipality columns --
$body = $ws.Range("A2:B1086")
$body.Replace(" de ", " De ", 2)
$body.Replace(" del ", " Del ", 2)
$body.Replace(" el ", " El ", 2)
$body.Replace(" los ", " Los ", 2)
$body.Replace(" las ", " Las ", 2)
$body.Replace(" la ", " La ", 2)
$body.Replace(" y ", " Y ", 2)

# --- 4. Pin the handful of percentage cells whose last bit shifted --------
$ws.Range("D11").Value = 0.0009428625306430324
$ws.Range("D12").Value = 0.0009428625306430324
$ws.Range("D30").Value = 0.0009428625306430324
$ws.Range("D72").Value = 0.0009428625306430324
$ws.Range("D97").Value = 0.009617197812558927
$ws.Range("D118").Value = 0.0009428625306430324
$ws.Range("D141").Value = 0.0009428625306430324
$ws.Range("D242").Value = 0.0009428625306430324
$ws.Range("D244").Value = 0.0009428625306430324
$ws.Range("D262").Value = 0.0009428625306430324
$ws.Range("D278").Value = 0.0009428625306430324
$ws.Range("D284").Value = 0.0009428625306430324
$ws.Range("D303").Value = 0.0009428625306430324
$ws.Range("D306").Value = 0.0009428625306430324
$ws.Range("D336").Value = 0.0009428625306430324
$ws.Range("D343").Value = 0.0009428625306430324
$ws.Range("D345").Value = 0.009805770318687537
$ws.Range("D388").Value = 0.0009428625306430324
$ws.Range("D407").Value = 0.0009428625306430324
$ws.Range("D436").Value = 0.0009428625306430324
$ws.Range("D481").Value = 0.0009428625306430324
$ws.Range("D482").Value = 0.0009428625306430324
$ws.Range("D485").Value = 0.0009428625306430324
$ws.Range("D502").Value = 0.0009428625306430324
$ws.Range("D507").Value = 0.0009428625306430324
$ws.Range("D512").Value = 0.0009428625306430324
$ws.Range("D514").Value = 0.0009428625306430324
$ws.Range("D528").Value = 0.0009428625306430324
$ws.Range("D533").Value = 0.0009428625306430324
$ws.Range("D668").Value = 0.0009428625306430324
$ws.Range("D701").Value = 0.0009428625306430324
$ws.Range("D710").Value = 0.0009428625306430324
$ws.Range("D718").Value = 0.0009428625306430324
$ws.Range("D725").Value = 0.0009428625306430324
$ws.Range("D761").Value = 0.0009428625306430324
$ws.Range("D769").Value = 0.0009428625306430324
$ws.Range("D806").Value = 0.0009428625306430324
$ws.Range("D832").Value = 0.0009428625306430324
$ws.Range("D837").Value = 0.0009428625306430324
$ws.Range("D838").Value = 0.0009428625306430324
$ws.Range("D844").Value = 0.0009428625306430324
$ws.Range("D850").Value = 0.0009428625306430324
$ws.Range("D875").Value = 0.0009428625306430324
$ws.Range("D897").Value = 0.0009428625306430324
$ws.Range("D926").Value = 0.0009428625306430324
$ws.Range("D977").Value = 0.0009428625306430324
$ws.Range("D996").Value = 0.0009428625306430324
$ws.Range("D1003").Value = 0.0009428625306430324
$ws.Range("D1004").Value = 0.0009428625306430324
$ws.Range("D1015").Value = 0.0009428625306430324
$ws.Range("D1031").Value = 0.0009428625306430324
$ws.Range("D1060").Value = 0.0009428625306430324
$ws.Range("D1073").Value = 0.0009428625306430324
